# daily auto push: 2026-01-26 13:52 UTC
#
# A new "2026/01/26" (Monday) data point (time 19, ranking 18) was logged
# alongside the workbook's existing "2026/01/26" entries. Insert it as a
# new row right after the last existing "2026/01/26" row (old row 702,
# the start of the "2026/12/29" block), shifting everything below it down
# by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 702 - this pushes the old rows 702:743
# down to 703:744 and extends the sheet's used range to row 744.
$ws.Rows.Item(702).Insert()

# Force column A to be stored as plain text so the date-like string isn't
# reinterpreted as a date serial number (matches every other cell in
# column A, which is stored as text).
$ws.Range("A702").NumberFormat = "@"
$ws.Range("A702").Value = "2026/01/26"
$ws.Range("B702").Value = "月"
$ws.Range("C702").Value = 19
$ws.Range("D702").Value = 18

# Drop the "Text" number-format style we applied above so the new row's
# formatting matches its neighbours (which carry no explicit style)
# while keeping the cell's stored content as text.
$ws.Range("A702").ClearFormats()
